# Auto-generated Excel COM-interop script
# Applies value updates to Sheets Sagittarius_Profits workbook per commit diff
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("J18").Value = 2333
$ws.Range("H18").Value = 1596.5
$ws.Range("N18").Value = -2901
$ws.Range("L18").Value = 2333
$ws.Range("M32").Value = -513.5
$ws.Range("J32").Value = 5198.6665
$ws.Range("L32").Value = 5198.6665
$ws.Range("N32").Value = -5850.6665
$ws.Range("K32").Value = 839.5
$ws.Range("H32").Value = 4108.875
$ws.Range("I32").Value = 839.5
$ws.Range("L86").Value = 501726
$ws.Range("H86").Value = 252605.38
$ws.Range("M86").Value = -2361.75
$ws.Range("J86").Value = 501726
$ws.Range("I86").Value = 3484.75
$ws.Range("K86").Value = 3484.75
$ws.Range("N86").Value = -503972
$ws.Range("I89").Value = 3484.75
$ws.Range("K89").Value = 17423.75
$ws.Range("J89").Value = 501726
$ws.Range("N89").Value = -2519862
$ws.Range("L89").Value = 2508630
$ws.Range("H89").Value = 252605.38
$ws.Range("M89").Value = -11807.75
$ws.Range("I113").Value = 4499
$ws.Range("K113").Value = 4499
$ws.Range("H113").Value = 4499
$ws.Range("M113").Value = -1245
$ws.Range("K116").Value = 16501.5
$ws.Range("M116").Value = -13059.5
$ws.Range("H116").Value = 7972.2856
$ws.Range("I116").Value = 16501.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H10").Value = 336283.16
$ws.Range("I10").Value = 502674.75
$ws.Range("K10").Value = 502674.75
$ws.Range("M10").Value = -502504.75
$ws.Range("M32").Value = -2502267.2
$ws.Range("K32").Value = 2502554.2
$ws.Range("H32").Value = 2657638.5
$ws.Range("I32").Value = 2502554.2
$ws.Range("M74").Value = -316.3
$ws.Range("K74").Value = 1190.3
$ws.Range("H74").Value = 1040.25
$ws.Range("I74").Value = 1190.3
$ws.Range("I77").Value = 1190.3
$ws.Range("K77").Value = 5951.5
$ws.Range("M77").Value = -1583.5
$ws.Range("H77").Value = 1040.25
$ws.Range("K110").Value = 2864
$ws.Range("H110").Value = 2864
$ws.Range("I110").Value = 2864
$ws.Range("M110").Value = -819

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("M80").Value = 827.6
$ws.Range("K80").Value = 170.4
$ws.Range("I80").Value = 170.4
$ws.Range("H80").Value = 331.7143
$ws.Range("I83").Value = 170.4
$ws.Range("K83").Value = 852
$ws.Range("M83").Value = 4140
$ws.Range("H83").Value = 331.7143
$ws.Range("K94").Value = 392.2
$ws.Range("N94").Value = -1725.5
$ws.Range("L94").Value = 823.5
$ws.Range("H94").Value = 464.08334
$ws.Range("I94").Value = 392.2
$ws.Range("M94").Value = 58.80000000000001
$ws.Range("J94").Value = 823.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H25").Value = 736.1667
$ws.Range("M25").Value = -562.1667
$ws.Range("K25").Value = 736.1667
$ws.Range("I25").Value = 736.1667
$ws.Range("M31").Value = -2863.25
$ws.Range("H31").Value = 3158.25
$ws.Range("I31").Value = 3158.25
$ws.Range("K31").Value = 3158.25
$ws.Range("I34").Value = 3158.25
$ws.Range("M34").Value = -2956.25
$ws.Range("H34").Value = 3158.25
$ws.Range("K34").Value = 3158.25
$ws.Range("M41").ClearContents()
$ws.Range("I41").Value = 0
$ws.Range("K41").Value = 0
$ws.Range("H41").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("J50").Value = 0
$ws.Range("N50").ClearContents()
$ws.Range("H50").Value = 0
$ws.Range("H62").Value = 0
$ws.Range("N62").ClearContents()
$ws.Range("L62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("N65").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("M99").Value = -384.4286
$ws.Range("K99").Value = 1882.4286
$ws.Range("N99").Value = -5053.1428
$ws.Range("J99").Value = 2057.1428
$ws.Range("L99").Value = 2057.1428
$ws.Range("H99").Value = 1969.7858
$ws.Range("I99").Value = 1882.4286
$ws.Range("I122").Value = 2000
$ws.Range("H122").Value = 2757.1428
$ws.Range("M122").Value = -3550
$ws.Range("K122").Value = 6000
$ws.Range("I126").Value = 1882.4286
$ws.Range("K126").Value = 5647.2858
$ws.Range("N126").Value = -11111.4284
$ws.Range("H126").Value = 1969.7858
$ws.Range("L126").Value = 6171.428400000001
$ws.Range("M126").Value = -3177.2858
$ws.Range("J126").Value = 2057.1428

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H76").Value = 16002.167
$ws.Range("H79").Value = 16002.167
$ws.Range("M111").Value = 1568.5
$ws.Range("I111").Value = 499.5
$ws.Range("K111").Value = 1498.5
$ws.Range("H111").Value = 499.5
$ws.Range("I119").Value = 3999
$ws.Range("M119").Value = -7159
$ws.Range("K119").Value = 11997
$ws.Range("H119").Value = 3999
$ws.Range("K120").Value = 8998.5
$ws.Range("H120").Value = 2999.5
$ws.Range("M120").Value = -4160.5
$ws.Range("I120").Value = 2999.5
$ws.Range("H136").Value = 2015
$ws.Range("I136").Value = 2015
$ws.Range("K136").Value = 6045
$ws.Range("M136").Value = -945

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("M132").Value = -24120.8
$ws.Range("I132").Value = 8883.6
$ws.Range("H132").Value = 8883.6
$ws.Range("K132").Value = 26650.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("K7").Value = 6899.3335
$ws.Range("N7").Value = -8087
$ws.Range("J7").Value = 7863
$ws.Range("L7").Value = 7863
$ws.Range("I7").Value = 6899.3335
$ws.Range("M7").Value = -6787.3335
$ws.Range("H7").Value = 7429.35
$ws.Range("H82").Value = 1232.3334
$ws.Range("I82").Value = 899
$ws.Range("K82").Value = 899
$ws.Range("M82").Value = -538
$ws.Range("K85").Value = 899
$ws.Range("I85").Value = 899
$ws.Range("H85").Value = 1232.3334
$ws.Range("M85").Value = 349
$ws.Range("I126").Value = 6899.3335
$ws.Range("K126").Value = 20698.0005
$ws.Range("N126").Value = -28529
$ws.Range("H126").Value = 7429.35
$ws.Range("L126").Value = 23589
$ws.Range("M126").Value = -18228.0005
$ws.Range("J126").Value = 7863
$ws.Range("H136").Value = 4156
$ws.Range("I136").Value = 4156
$ws.Range("K136").Value = 12468
$ws.Range("M136").Value = -9918

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I6").Value = 3247.5
$ws.Range("K6").Value = 3247.5
$ws.Range("M6").Value = -3132.5
$ws.Range("H6").Value = 3321.125
$ws.Range("I62").Value = 9161.857
$ws.Range("K62").Value = 9161.857
$ws.Range("H62").Value = 11913.2
$ws.Range("M62").Value = -8537.857
$ws.Range("I65").Value = 9161.857
$ws.Range("M65").Value = -42689.285
$ws.Range("K65").Value = 45809.285
$ws.Range("H65").Value = 11913.2
$ws.Range("M132").Value = -5449.25
$ws.Range("I132").Value = 2659.75
$ws.Range("H132").Value = 2659.75
$ws.Range("K132").Value = 7979.25
$ws.Range("H136").Value = 2147.4119
$ws.Range("I136").Value = 2162.875
$ws.Range("K136").Value = 6488.625
$ws.Range("M136").Value = -3938.625
